$wb = $excel.ActiveWorkbook

# The sheet currently named "InvestmentOpn" becomes "Investment_Opening_Cash".
$ws = $wb.Worksheets.Item("InvestmentOpn")
$ws.Name = "Investment_Opening_Cash"

# Update the form-name cell (A2) to match the new module name.
$ws.Range("A2").Value = "Investment_Opening_Cash"

# Row 2 got taller in the new layout.
$ws.Rows.Item(2).RowHeight = 60

# This sheet becomes the active tab, with the selection moved to I11.
$ws.Activate()
$ws.Range("I11").Select()
